$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 9868.177
$ws.Range("I28").Value = 706.625
$ws.Range("J28").Value = 18011.777
$ws.Range("K28").Value = 706.625
$ws.Range("L28").Value = 18011.777
$ws.Range("M28").Value = -221.625
$ws.Range("N28").Value = -18981.777

# Row 33
$ws.Range("H33").Value = 439.05264
$ws.Range("I33").Value = 464.82352
$ws.Range("K33").Value = 464.82352
$ws.Range("M33").Value = -235.82352

# Row 51
$ws.Range("H51").Value = 1694.8334
$ws.Range("I51").Value = 1694.8334
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 1694.8334
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -1210.8334
$ws.Range("N51").Value = $null   # was -2468

# Row 98
$ws.Range("H98").Value = 878.9259
$ws.Range("I98").Value = 706.35
$ws.Range("J98").Value = 1372
$ws.Range("K98").Value = 706.35
$ws.Range("L98").Value = 1372
$ws.Range("M98").Value = 791.65
$ws.Range("N98").Value = -4368

# Row 122
$ws.Range("H122").Value = 878.9259
$ws.Range("I122").Value = 706.35
$ws.Range("J122").Value = 1372
$ws.Range("K122").Value = 2119.05
$ws.Range("L122").Value = 4116
$ws.Range("M122").Value = 330.9499999999998
$ws.Range("N122").Value = -9016

# Row 125
$ws.Range("H125").Value = 2134.8572
$ws.Range("I125").Value = 882
$ws.Range("J125").Value = 2636
$ws.Range("K125").Value = 7938
$ws.Range("L125").Value = 23724
$ws.Range("M125").Value = -5478
$ws.Range("N125").Value = -28644

# Row 127
$ws.Range("H127").Value = 1487.9565
$ws.Range("I127").Value = 369.4
$ws.Range("J127").Value = 1798.6666
$ws.Range("K127").Value = 1108.2
$ws.Range("L127").Value = 5395.9998
$ws.Range("M127").Value = 3851.8
$ws.Range("N127").Value = -15315.9998

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1262.5
$ws.Range("I2").Value = 1262.5
$ws.Range("K2").Value = 1262.5
$ws.Range("M2").Value = -1149.5

# Row 24
$ws.Range("H24").Value = 54375
$ws.Range("J24").Value = 54375
$ws.Range("L24").Value = 54375
$ws.Range("N24").Value = -55123

# Row 32
$ws.Range("H32").Value = 1414273.5
$ws.Range("I32").Value = 1561555.5
$ws.Range("J32").Value = 59278.8
$ws.Range("K32").Value = 1561555.5
$ws.Range("L32").Value = 59278.8
$ws.Range("M32").Value = -1561268.5
$ws.Range("N32").Value = -59852.8

# Row 45
$ws.Range("H45").Value = 1893.091
$ws.Range("I45").Value = 2279.2
$ws.Range("J45").Value = 1571.3334
$ws.Range("K45").Value = 2279.2
$ws.Range("L45").Value = 1571.3334
$ws.Range("M45").Value = -1902.2
$ws.Range("N45").Value = -2325.3334

# Row 63
$ws.Range("H63").Value = 7111.421
$ws.Range("I63").Value = 5595.6665
$ws.Range("J63").Value = 7811
$ws.Range("K63").Value = 5595.6665
$ws.Range("L63").Value = 7811
$ws.Range("M63").Value = -4909.6665
$ws.Range("N63").Value = -9183

# Row 66
$ws.Range("H66").Value = 7111.421
$ws.Range("I66").Value = 5595.6665
$ws.Range("J66").Value = 7811
$ws.Range("K66").Value = 27978.3325
$ws.Range("L66").Value = 39055
$ws.Range("M66").Value = -24546.3325
$ws.Range("N66").Value = -45919

# Row 97
$ws.Range("H97").Value = 1627.6923
$ws.Range("I97").Value = 1495
$ws.Range("J97").Value = 1840
$ws.Range("K97").Value = 1495
$ws.Range("L97").Value = 1840
$ws.Range("M97").Value = -999
$ws.Range("N97").Value = -2832

# Row 100
$ws.Range("H100").Value = 54375
$ws.Range("J100").Value = 54375
$ws.Range("L100").Value = 54375
$ws.Range("N100").Value = -56539

# Row 116
$ws.Range("H116").Value = 1262.5
$ws.Range("I116").Value = 1262.5
$ws.Range("K116").Value = 1262.5
$ws.Range("M116").Value = 1031.5

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1262.5
$ws.Range("I3").Value = 1262.5
$ws.Range("K3").Value = 1262.5
$ws.Range("M3").Value = -1148.5

# Row 82
$ws.Range("H82").Value = 16884.875
$ws.Range("I82").Value = 12371.333
$ws.Range("J82").Value = 19593
$ws.Range("K82").Value = 12371.333
$ws.Range("L82").Value = 19593
$ws.Range("M82").Value = -11988.333
$ws.Range("N82").Value = -20359

# Row 85
$ws.Range("H85").Value = 16884.875
$ws.Range("I85").Value = 12371.333
$ws.Range("J85").Value = 19593
$ws.Range("K85").Value = 12371.333
$ws.Range("L85").Value = 19593
$ws.Range("M85").Value = -11045.333
$ws.Range("N85").Value = -22245

# Row 94
$ws.Range("H94").Value = 908.11536
$ws.Range("I94").Value = 742.7368
$ws.Range("J94").Value = 1357
$ws.Range("K94").Value = 742.7368
$ws.Range("L94").Value = 1357
$ws.Range("M94").Value = -291.7368
$ws.Range("N94").Value = -2259

# Row 107
$ws.Range("H107").Value = 51479.3
$ws.Range("I107").Value = 77987.92
$ws.Range("K107").Value = 77987.92
$ws.Range("M107").Value = -76067.92

$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("H99").Value = 1780.303
$ws.Range("I99").Value = 1055.875
$ws.Range("J99").Value = 2012.12
$ws.Range("K99").Value = 1055.875
$ws.Range("L99").Value = 2012.12
$ws.Range("M99").Value = 442.125
$ws.Range("N99").Value = -5008.12

# Row 126
$ws.Range("H126").Value = 1780.303
$ws.Range("I126").Value = 1055.875
$ws.Range("J126").Value = 2012.12
$ws.Range("K126").Value = 3167.625
$ws.Range("L126").Value = 6036.36
$ws.Range("M126").Value = -697.625
$ws.Range("N126").Value = -10976.36

# Row 135
$ws.Range("H135").Value = 49234.75
$ws.Range("J135").Value = 49234.75
$ws.Range("L135").Value = 49234.75
$ws.Range("N135").Value = -59374.75

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1332.0769
$ws.Range("I5").Value = 574
$ws.Range("J5").Value = 2216.5
$ws.Range("K5").Value = 1722
$ws.Range("L5").Value = 6649.5
$ws.Range("M5").Value = -1610
$ws.Range("N5").Value = -6873.5

# Row 12
$ws.Range("H12").Value = 54.52381
$ws.Range("J12").Value = 59.666668
$ws.Range("L12").Value = 179.000004
$ws.Range("N12").Value = -525.000004

# Row 63
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").Value = $null   # was -3385

# Row 64
$ws.Range("H64").Value = 2233.3215
$ws.Range("J64").Value = 2389.12
$ws.Range("L64").Value = 7167.36
$ws.Range("N64").Value = -7707.36

# Row 66
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").Value = $null   # was -8658

# Row 67
$ws.Range("H67").Value = 2233.3215
$ws.Range("J67").Value = 2389.12
$ws.Range("L67").Value = 7167.36
$ws.Range("N67").Value = -9039.360000000001

# Row 69
$ws.Range("H69").Value = 2726
$ws.Range("I69").Value = 1012
$ws.Range("J69").Value = 3011.6667
$ws.Range("K69").Value = 3036
$ws.Range("L69").Value = 9035.000100000001
$ws.Range("M69").Value = -2225
$ws.Range("N69").Value = -10657.0001

# Row 70
$ws.Range("H70").Value = 6553.25
$ws.Range("I70").Value = 7404.3335
$ws.Range("K70").Value = 22213.0005
$ws.Range("M70").Value = -21898.0005

# Row 72
$ws.Range("H72").Value = 2726
$ws.Range("I72").Value = 1012
$ws.Range("J72").Value = 3011.6667
$ws.Range("K72").Value = 9108
$ws.Range("L72").Value = 27105.0003
$ws.Range("M72").Value = -5052
$ws.Range("N72").Value = -35217.0003

# Row 73
$ws.Range("H73").Value = 6553.25
$ws.Range("I73").Value = 7404.3335
$ws.Range("K73").Value = 22213.0005
$ws.Range("M73").Value = -21121.0005

# Row 98
$ws.Range("H98").Value = 833.5
$ws.Range("J98").Value = 866
$ws.Range("L98").Value = 2598
$ws.Range("N98").Value = -5594

# Row 122
$ws.Range("H122").Value = 2946.1592
$ws.Range("I122").Value = 462.75
$ws.Range("J122").Value = 4365.25
$ws.Range("K122").Value = 4164.75
$ws.Range("L122").Value = 39287.25
$ws.Range("M122").Value = -1714.75
$ws.Range("N122").Value = -44187.25

# Row 131
$ws.Range("H131").Value = 4476.5454
$ws.Range("I131").Value = 257.27274
$ws.Range("J131").Value = 6586.1816
$ws.Range("K131").Value = 771.81822
$ws.Range("L131").Value = 19758.5448
$ws.Range("M131").Value = 4268.18178
$ws.Range("N131").Value = -29838.5448

# Row 135
$ws.Range("H135").Value = 1332.0769
$ws.Range("I135").Value = 574
$ws.Range("J135").Value = 2216.5
$ws.Range("K135").Value = 5166
$ws.Range("L135").Value = 19948.5
$ws.Range("M135").Value = -2631
$ws.Range("N135").Value = -25018.5

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1758.3636
$ws.Range("I102").Value = 1390
$ws.Range("J102").Value = 2290.4443
$ws.Range("K102").Value = 1390
$ws.Range("L102").Value = 2290.4443
$ws.Range("M102").Value = 232
$ws.Range("N102").Value = -5534.4443

# Row 122
$ws.Range("H122").Value = 1581.5454
$ws.Range("I122").Value = 1581.5454
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4744.6362
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2294.6362
$ws.Range("N122").Value = $null   # was -10597

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 9234.556
$ws.Range("I7").Value = 3558
$ws.Range("K7").Value = 3558
$ws.Range("M7").Value = -3446

# Row 40
$ws.Range("H40").Value = 4214.857
$ws.Range("I40").Value = 4401
$ws.Range("J40").Value = 3966.6667
$ws.Range("K40").Value = 4401
$ws.Range("L40").Value = 3966.6667
$ws.Range("M40").Value = -4265
$ws.Range("N40").Value = -4238.6667

# Row 122
$ws.Range("H122").Value = 2938
$ws.Range("I122").Value = 2876
$ws.Range("K122").Value = 8628
$ws.Range("M122").Value = -6178

# Row 126
$ws.Range("H126").Value = 9234.556
$ws.Range("I126").Value = 3558
$ws.Range("K126").Value = 10674
$ws.Range("M126").Value = -8204

# Row 136
$ws.Range("H136").Value = 6947016.5
$ws.Range("I136").Value = 2262.2222
$ws.Range("J136").Value = 27781278
$ws.Range("K136").Value = 6786.6666
$ws.Range("L136").Value = 83343834
$ws.Range("M136").Value = -4236.6666
$ws.Range("N136").Value = -83348934

$ws = $wb.Worksheets.Item("WVR")
# Row 15
$ws.Range("H15").Value = 31251.75
$ws.Range("J15").Value = 31251.75
$ws.Range("L15").Value = 31251.75
$ws.Range("N15").Value = -31827.75

# Row 122
$ws.Range("H122").Value = 2819.4119
$ws.Range("I122").Value = 2910
$ws.Range("J122").Value = 2525
$ws.Range("K122").Value = 8730
$ws.Range("L122").Value = 7575
$ws.Range("M122").Value = -6280
$ws.Range("N122").Value = -12475

# Row 126
$ws.Range("H126").Value = 1130.8
$ws.Range("I126").Value = 1068
$ws.Range("J126").Value = 1225
$ws.Range("K126").Value = 3204
$ws.Range("L126").Value = 3675
$ws.Range("M126").Value = -734
$ws.Range("N126").Value = -8615
